$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Forces the cell to keep an exact text representation (preserving leading/
    # trailing zeros, thousand-dot separators, etc.) instead of Excel auto-coercing
    # the string into a number, and then clears the temporary formatting so the
    # cell keeps its original (unstyled) look.
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '43.221.46'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '2.353.33'
$ws.Range('E3').Value = '  +5.09%  '
$ws.Range('E4').Value = '  +0.07%  '
Set-TextValue 'D5' '234.52'
$ws.Range('E5').Value = '  +1.98%  '
Set-TextValue 'D6' '0.650'
$ws.Range('E6').Value = '  +1.75%  '
Set-TextValue 'D7' '71.50'
$ws.Range('E7').Value = '  +14.10%  '
Set-TextValue 'D9' '0.497'
$ws.Range('E9').Value = '  +13.08%  '
Set-TextValue 'D10' '0.0973'
$ws.Range('E10').Value = '  +2.03%  '
Set-TextValue 'D11' '27.28'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('E12').Value = '  +2.58%  '
$ws.Range('D13').Value = '2.705.78'
$ws.Range('E13').Value = '  +5.33%  '
Set-TextValue 'D14' '16.26'
$ws.Range('E14').Value = '  +5.90%  '
Set-TextValue 'D15' '6.35'
$ws.Range('E15').Value = '  +5.63%  '
Set-TextValue 'D16' '0.865'
$ws.Range('E16').Value = '  +5.02%  '
$ws.Range('D17').Value = '2.368.88'
$ws.Range('E17').Value = '  +5.85%  '
$ws.Range('D18').Value = '43.270.59'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('E19').Value = '  +4.83%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D20' '6.35'
$ws.Range('E20').Value = '  +4.86%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D21' '74.69'
$ws.Range('E21').Value = '  +2.78%  '
Set-TextValue 'D22' '250.27'
$ws.Range('E22').Value = '  +2.10%  '
$ws.Range('B23').Value = 'WEMIXToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D23' '3.81'
$ws.Range('E23').Value = '  +3.88%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D24' '1.00'
$ws.Range('E24').Value = '  +0.04%  '
Set-TextValue 'D25' '2.46'
$ws.Range('E25').Value = '  +2.56%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D26' '2.25'
$ws.Range('E26').Value = '  -1.06%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D27' '10.05'
$ws.Range('E27').Value = '  +3.97%  '
Set-TextValue 'D28' '22.39'
$ws.Range('E28').Value = '  +4.40%  '
Set-TextValue 'D29' '172.39'
$ws.Range('E29').Value = '  +0.48%  '
Set-TextValue 'D30' '1.54'
$ws.Range('E30').Value = '  +10.15%  '
$ws.Range('E31').Value = '  +1.72%  '
$ws.Range('E32').Value = '  +2.83%  '
Set-TextValue 'D33' '5.00'
$ws.Range('E33').Value = '  +3.36%  '
Set-TextValue 'D34' '0.0691'
$ws.Range('E34').Value = '  +3.25%  '
$ws.Range('E35').Value = '  +4.55%  '
Set-TextValue 'D36' '3.73'
$ws.Range('E36').Value = '  +3.87%  '
Set-TextValue 'D37' '6.54'
$ws.Range('E37').Value = '  +4.35%  '
Set-TextValue 'D38' '2.43'
$ws.Range('E38').Value = '  +7.45%  '
Set-TextValue 'D39' '0.0255'
$ws.Range('E39').Value = '  +2.39%  '
Set-TextValue 'D40' '19.01'
$ws.Range('E40').Value = '  +13.10%  '
$ws.Range('B41').Value = 'BinanceUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D41' '1.00'
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D42' '8.92'
$ws.Range('E42').Value = '  +4.91%  '
$ws.Range('E43').Value = '  +0.27%  '
Set-TextValue 'D44' '99.16'
$ws.Range('E44').Value = '  +3.52%  '
$ws.Range('E45').Value = '  +9.67%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D46' '0.0961'
$ws.Range('E46').Value = '  +2.94%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D47' '1.21'
$ws.Range('E47').Value = '  +3.52%  '
$ws.Range('D48').Value = '1.443.67'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').Value = '2.578.99'
$ws.Range('E49').Value = '  +5.51%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D50' '2.77'
$ws.Range('E50').Value = '  +0.83%  '
$ws.Range('B51').Value = 'TerraClassic'
$ws.Range('C51').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextValue 'D51' '0.000202'
$ws.Range('E51').Value = '  -3.36%  '
